$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 30 data rows (rows 2-31) down by 11 rows to rows 13-42,
# preserving values + formatting in a single copy/paste operation (keeps
# the workbook's existing style table intact - no new styles minted).
$ws.Range("A2:E31").Copy($ws.Range("A13"))

# Fill the newly freed rows 2-12 with the additional (earlier / backward
# extension) forecast-vector data.
$newRows = @(
    @(31228, 1985, 1.777145513200118, 1986, 1.897195788720762),
    @(31593, 1986, 1.520629350269975, 1987, 1.941216181406236),
    @(31958, 1987, 0.006150708382479664, 1988, 2.273434064658209),
    @(32324, 1988, 2.284628917872622, 1989, 2.270390965314983),
    @(32689, 1989, 3.661580277249166, 1990, 2.193258610001214),
    @(33054, 1990, 4.130619852766437, 1991, 2.272904440822465),
    @(33419, 1991, 6.364491101711689, 1992, 2.793798186209284),
    @(33785, 1992, 2.932796654414149, 1993, 2.584450468619459),
    @(34150, 1993, -1.06363680093724, 1994, 2.368493192930488),
    @(34515, 1994, 2.479893153134016, 1995, 2.567096653116252),
    @(34880, 1995, 2.432437183852798, 1996, 2.872765583543457)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 2 + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
